# Auto-generated update of crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '307.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-2.64%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '40.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.16%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.038'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.40%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07638'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-3.52%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.230'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-3.42%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.596'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-5.19%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9111'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.60%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.452'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-4.36%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1021'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-8.06%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1766'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.64%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09066'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.96%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04381'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.45%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1055'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.07%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001248'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.24%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005823'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.62%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.371'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.52%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.41%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.748'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-6.85%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1355'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.70%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2715'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.91%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04150'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.39%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.63%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004025'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-3.62%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '5.82%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0002998'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.96%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02424'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-1.67%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05179'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-2.95%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007788'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.10%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1310'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.38%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007048'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-8.31%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001947'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.12%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007748'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-6.08%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3061'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.57%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006362'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-7.32%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000748'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-1.08%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.005342'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '57.44%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004387'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '5.29%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002094'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-1.08%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001994'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-1.08%'
